$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRA form (empty)")

$ws.Range("B10").Value = "Working with sound >85 dB"
$ws.Range("C10").Value = "Hearing loss"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 6
$ws.Range("I10").Value = "Wear ear protection above 85 dB"
$ws.Range("H10").Value = "Monitor sound level, avoid playing sound above 85 dB if not necessary, do not play sound above 100 dB"
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 6
$ws.Range("L10").Value = 0.5
$ws.Range("N10").Value = "Immediately see a doctor"
$ws.Range("O10").Value = "Lotte, Max, Sam, Thijs"

$ws.Rows.Item(10).EntireRow.AutoFit() | Out-Null
